$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = 44181
$ws.Range("M3").Value = 30

# Row 4
$ws.Range("D4").Value = 44175
$ws.Range("M4").Value = 25
$ws.Range("N4").Value = 20000
$ws.Range("O4").Value = 20000
$ws.Range("P4").Value = 20000
$ws.Range("S4").Value = 4000

# Row 5
$ws.Range("D5").Value = 44186
$ws.Range("M5").Value = 40
$ws.Range("N5").Value = 15000
$ws.Range("O5").Value = 15000
$ws.Range("P5").Value = 15000
$ws.Range("S5").Value = 3000

# Row 6
$ws.Range("D6").Value = 44188
$ws.Range("M6").Value = 30

# Row 7
$ws.Range("D7").Value = 44189
$ws.Range("M7").Value = 40
$ws.Range("N7").Value = 15000
$ws.Range("O7").Value = 15000
$ws.Range("P7").Value = 15000
$ws.Range("S7").Value = 3000

# Row 9
$ws.Range("D9").Value = 44179
$ws.Range("M9").Value = 45
$ws.Range("N9").Value = 20000
$ws.Range("O9").Value = 20000
$ws.Range("P9").Value = 20000
$ws.Range("S9").Value = 4000
